$d = $word.ActiveDocument

# 1. Objective paragraph: rewrite the objective sentence.
$d.Content.Find.Execute(
    "Seeking to contribute technical and creative problem-solving skills to NPR" + [char]0x2019 + "s innovative work in multiple fields including digital media and technology. Highly motivated and detail-oriented Software Engineering student with hands-on experience in web development looking to apply learned skills to make an impact.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Seeking to contribute technical and creative problem-solving skills to innovative work in technology. Highly motivated and detail-oriented Software Engineering student with hands-on experience in software development looking to apply learned skills to make an impact.",
    2) | Out-Null

# 2. Education: drop the " | Expected Graduation: June 2026" suffix after the university name.
$d.Content.Find.Execute(
    "Western Governors University | Expected Graduation: June 2026",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Western Governors University",
    2) | Out-Null

# 3. Relevant Coursework: update the course list.
$d.Content.Find.Execute(
    "Relevant Coursework: Front-End Web Development, User Interface Design, User Experience Design, JavaScript Programming, Cloud Foundations, Software Design and Quality Assurance, Version Control",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Relevant Coursework: Data Structures & Algorithms, Front-End Web Development, JavaScript Programming, Intro to Programming in Python, Cloud Foundations, Advanced Data Management, Software Design and Quality Assurance, Version Control",
    2) | Out-Null

# 4. Remove the "Other: Microsoft Office, Excel, Zoom, Webex" bullet entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Other: Microsoft Office*") {
        $p.Range.Delete()
        break
    }
}

# 5. Certificates: add a new "CompTIA Project+" bullet after the ITIL Foundation entry.
$idx = 0
$target = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "ITIL Foundation in IT Service Management*") {
        $target = $idx
    }
}
if ($target -ge 0) {
    $p = $d.Paragraphs.Item($target)
    $p.Range.InsertParagraphAfter() | Out-Null
    $d.Paragraphs.Item($target + 1).Range.Text = "CompTIA Project+"
}
